$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.519.31"
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.650.31"
$ws.Range("E3").Value = "  +1.79%  "

$ws.Range("E4").Value = "  +0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "299.34"
$ws.Range("E6").Value = "  -1.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3792"
$ws.Range("E7").Value = "  +1.23%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.77"
$ws.Range("E8").Value = "  -2.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3548"
$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08091"
$ws.Range("E10").Value = "  -0.62%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.220"
$ws.Range("E11").Value = "  -0.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.99"
$ws.Range("E13").Value = "  -0.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.376"
$ws.Range("E14").Value = "  -1.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.404"
$ws.Range("E15").Value = "  +1.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001198"
$ws.Range("E16").Value = "  -1.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.650.29"
$ws.Range("E17").Value = "  +2.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.25"
$ws.Range("E18").Value = "  +2.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06988"
$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.792"
$ws.Range("E20").Value = "  +1.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.47"
$ws.Range("E21").Value = "  +0.68%  "

$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.61"
$ws.Range("E23").Value = "  +1.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.531.24"
$ws.Range("E24").Value = "  +1.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.484"
$ws.Range("E25").Value = "  -1.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.938"
$ws.Range("E26").Value = "  -4.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.02"
$ws.Range("E27").Value = "  +0.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.55"
$ws.Range("E28").Value = "  -0.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.232"
$ws.Range("E29").Value = "  +1.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.95"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.836.69"
$ws.Range("E31").Value = "  +2.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.939"
$ws.Range("E32").Value = "  +6.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.172"
$ws.Range("E33").Value = "  +7.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.75"
$ws.Range("E34").Value = "  +0.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.033"
$ws.Range("E35").Value = "  -5.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02731"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08718"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2447"
$ws.Range("E38").Value = "  -0.69%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.960"
$ws.Range("E39").Value = "  +1.42%  "

$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "13.12"
$ws.Range("E40").Value = "  +4.76%  "

$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06845"
$ws.Range("E41").Value = "  -1.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6911"
$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.316"
$ws.Range("E43").Value = "  -0.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.75"
$ws.Range("E44").Value = "  +1.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6423"
$ws.Range("E45").Value = "  +1.02%  "

$ws.Range("E46").Value = "  +0.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.260"
$ws.Range("E47").Value = "  -0.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.919"
$ws.Range("E48").Value = "  -0.56%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07875"
$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.72"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.176"
$ws.Range("E51").Value = "  +0.61%  "
